$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: the "Validate the login functionality ... FWC software tool"
# table-cell paragraph is rewritten. The old paragraph held the _GoBack
# bookmark straddling "of  " / "the"; the new paragraph drops that bookmark,
# splits the original sentence into several runs (Word's grammar checker
# wraps "of  the" in <w:proofErr> start/end markers) and appends a second
# sentence about valid credentials (also proofErr-wrapped at the end).
# ---------------------------------------------------------------------------
$oldText = "Validate the login functionality of  the FWC software tool"
$cell = $d.Content
$found = $cell.Find.Execute($oldText, $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the login-functionality sentence to rewrite."
}

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$paraXml = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00175E0B" w:rsidRDefault="00175E0B" w:rsidP="00175E0B">' + `
  '<w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Validate the login functionality </w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:t>of  the</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> FWC software tool</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">. Any registered user should be able to login with the valid username, </w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:t>password .</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '</w:p>'

[void]$cell.InsertXML($paraXml)

# ---------------------------------------------------------------------------
# Change 2: the _GoBack bookmark (dropped from the cell above) reappears on
# its own in the last, otherwise-empty paragraph at the very end of the
# document body. Bookmarks.Add ignores genuinely empty ranges in this
# runtime, so stage a one-character placeholder, bookmark it, then delete
# the placeholder text (the bookmark start/end tags survive the deletion).
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$lastRange = $d.Paragraphs.Item($n).Range
$lastRange.InsertBefore("X")
$placeholder = $d.Range($lastRange.Start, $lastRange.Start + 1)
[void]$d.Bookmarks.Add("_GoBack", $placeholder)
$placeholder = $d.Range($lastRange.Start, $lastRange.Start + 1)
$placeholder.Text = ""
